$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6233027577400208
$ws.Range("B1").Value = 3.548562526702881
$ws.Range("C1").Value = 2.697499990463257
$ws.Range("D1").Value = 0.9549130201339722
$ws.Range("E1").Value = 0.9180357456207275
